$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H column values (subcategory) for rows where "line graph" -> "line graph(s)"
$lineGraphRows = @(3, 4, 7, 8, 9, 48)
foreach ($r in $lineGraphRows) {
    $ws.Cells.Item($r, 8).Value = "line graph(s)"
}

# Update H28: "data structure" -> "data display"
$ws.Cells.Item(28, 8).Value = "data display"

# Remove the "is_viewed" column (column I) entirely
$ws.Columns.Item(9).Delete()
